# Adds documentation rows for the two new NetMHCpan "Of" (offset) columns:
#   - NetMHCpan_bestRank_Of      (inserted right after NetMHCpan_bestRank_Icore)
#   - NetMHCpan_bestAffinity_Of  (inserted right after NetMHCpan_bestAffinity_Icore)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Neoantigen")

$offsetDescription = "Starting position offset of the core in the peptide (0 based) "

# --- Insert row for NetMHCpan_bestAffinity_Of right before the current row 23
# (which currently holds NetMHCpan_bestAffinity_Gp). This shifts row 23 and below down by one.
$ws.Rows.Item(23).Insert()
$ws.Cells.Item(23, 1).Value = "NetMHCpan_bestAffinity_Of"
$ws.Cells.Item(23, 2).Value = $offsetDescription
$ws.Cells.Item(23, 3).Value = " MHC I binding with netMHCpan    "

# --- Insert row for NetMHCpan_bestRank_Of right before row 16
# (which still holds NetMHCpan_bestRank_Gp, unaffected by the insert above since 16 < 23).
# This shifts row 16 and below (including the new row we just added) down by one.
$ws.Rows.Item(16).Insert()
$ws.Cells.Item(16, 1).Value = "NetMHCpan_bestRank_Of"
$ws.Cells.Item(16, 2).Value = $offsetDescription
$ws.Cells.Item(16, 3).Value = " MHC I binding with netMHCpan    "

# --- Update the view so it matches the scrolled/selected position after the edit.
$ws.Application.ActiveWindow.ScrollRow = 34
$ws.Range("B50").Select()
